$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$titles = @(
    "Gigantyczna lawina w rejonie Morskiego Oka. Jest nagranie",
    "Ulubiony klub jachtowy Rosjan w płomieniach",
    "Problemy w Willingen! Polacy poradzą sobie w kwalifikacjach?",
    "Media: Reznikow straci stanowisko. Podano nazwisko następcy",
    "Depardieu: Dla mnie nic się nie zmieniło. Nadal jestem Rosjaninem",
    "Politycy odpowiadają premierowi. Czarzasty: Małpie brzytwy się nie daje",
    "Patrioty zmierzają do Warszawy. Zostaną rozstawione na lotnisku",
    "Te osoby mogą podwyższyć sobie emeryturę. Wystarczy jeden wniosek",
    "Filmowe tytuły, z których śmieją się wszyscy. Kto je wybiera?"
)

for ($i = 0; $i -lt $titles.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $titles[$i]
}
